# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" sheets to reflect the latest scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 521   # 南宁·小蜜蜂动漫嘉年华2.0: 520 -> 521
$ws1.Range("F5").Value = 218   # 南宁·漫控嘉年华09...: 215 -> 218
$ws1.Range("F7").Value = 228   # 广西·首届明日方舟only展: 227 -> 228
$ws1.Range("F8").Value = 2230  # 南宁·AB动漫游戏嘉年华: 2227 -> 2230
$ws1.Range("F10").Value = 5495 # 南宁·第十九届良牙动漫夏季盛典: 5480 -> 5495
$ws1.Range("F12").Value = 362  # 南宁·蔚蓝档案only: 361 -> 362

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 521   # 南宁·小蜜蜂动漫嘉年华2.0: 520 -> 521
$ws4.Range("F6").Value = 218   # 南宁·漫控嘉年华09...: 215 -> 218
$ws4.Range("F8").Value = 228   # 广西·首届明日方舟only展: 227 -> 228
$ws4.Range("F11").Value = 2230 # 南宁·AB动漫游戏嘉年华: 2227 -> 2230
$ws4.Range("F13").Value = 5495 # 南宁·第十九届良牙动漫夏季盛典: 5480 -> 5495
$ws4.Range("F15").Value = 362  # 南宁·蔚蓝档案only: 361 -> 362

$wb.Save()
